$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.124112129211426
$ws.Range("B1").Value = 2.503989458084106
$ws.Range("C1").Value = 6.106705188751221
$ws.Range("D1").Value = 2.164223432540894
$ws.Range("E1").Value = 1.246336460113525
